# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund-holdings detail, same shape as the
# other quarterly sheets) right before the "总计" (totals) sheet, and adds
# a corresponding summary row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Reuse the existing bold/bordered header style (as used by the other
# quarterly sheets, e.g. "2021-Q4") instead of minting a new one: copy
# just the formatting (no values) from that sheet onto the new sheet.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A6").Copy()
$q1.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Data rows. Fund code / size / position figures are kept as plain
#    text (a leading "'" forces text entry without altering any cell
#    style) so codes like "005197" and figures like "0.10" keep their
#    original, exact textual representation instead of being coerced
#    into numbers.
# ---------------------------------------------------------------------
$rows = @(
    @("005197", "工银瑞信沪港深精选灵活配置混合A",     "7.18", "93.54", "4.07", "0.2922", 7),
    @("005504", "汇添富沪港深大盘价值混合",             "3.49", "92.33", "5.14", "0.1794", 7),
    @("015119", "汇添富沪港深大盘价值混合D",            "3.49", "92.33", "5.14", "0.1794", 7),
    @("005198", "工银瑞信沪港深精选灵活配置混合C",     "2.66", "93.54", "4.07", "0.1083", 7),
    @("160125", "南方香港优选股票QDII-LOF",             "2.46", "91.14", "3.61", "0.0888", 5),
    @("161229", "国投瑞银中国价值发现股票QDII-LOF",     "1.47", "92.83", "5.27", "0.0775", 5),
    @("004532", "民生加银中证港股通高股息精选指数A",   "0.26", "94.88", "3.53", "0.0092", 9),
    @("004533", "民生加银中证港股通高股息精选指数C",   "0.10", "94.88", "3.53", "0.0035", 9),
    @("005770", "信达澳银中证沪港深高股息精选指数",     "0.01", "92.47", "2.79", "0.0003", 2)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rec = $rows[$i]

    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = "'" + $rec[0]
    $q1.Cells.Item($r, 3).Value = $rec[1]
    $q1.Cells.Item($r, 4).Value = "'" + $rec[2]
    $q1.Cells.Item($r, 5).Value = "'" + $rec[3]
    $q1.Cells.Item($r, 6).Value = "'" + $rec[4]
    $q1.Cells.Item($r, 7).Value = "'" + $rec[5]
    $q1.Cells.Item($r, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 4. Prepend the new quarter to the "总计" (totals) summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.9399999999999999

# The index column (A) for the quarters that shifted down by one row
# needs to be bumped so it keeps counting up from 0.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
